{"js": "// The document consists of a centered date paragraph followed by a 5x? table\n// of two-digit division problems. Every non-empty paragraph's text changes,\n// in document order, from an \"old\" value to a \"new\" value (1 date + 25 table\n// cells). Walk body.paragraphs (which includes paragraphs nested in table\n// cells, in document order) and rewrite each non-empty one in turn.\nconst pairs = [[\"2023-08-02 Wednesday\", \"2023-08-03 Thursday\"], [\"89\u00f74=22, 1\", \"87\u00f73=29, 0\"], [\"89\u00f77=12, 5\", \"12\u00f78=1, 4\"], [\"36\u00f72=18, 0\", \"97\u00f72=48, 1\"], [\"28\u00f75=5, 3\", \"72\u00f79=8, 0\"], [\"86\u00f73=28, 2\", \"38\u00f73=12, 2\"], [\"99\u00f79=11, 0\", \"16\u00f76=2, 4\"], [\"38\u00f75=7, 3\", \"62\u00f75=12, 2\"], [\"98\u00f78=12, 2\", \"60\u00f75=12, 0\"], [\"82\u00f79=9, 1\", \"71\u00f74=17, 3\"], [\"80\u00f75=16, 0\", \"12\u00f75=2, 2\"], [\"83\u00f74=20, 3\", \"27\u00f72=13, 1\"], [\"79\u00f72=39, 1\", \"89\u00f75=17, 4\"], [\"18\u00f73=6, 0\", \"51\u00f79=5, 6\"], [\"89\u00f77=12, 5\", \"75\u00f76=12, 3\"], [\"93\u00f73=31, 0\", \"24\u00f75=4, 4\"], [\"64\u00f79=7, 1\", \"93\u00f74=23, 1\"], [\"83\u00f74=20, 3\", \"32\u00f77=4, 4\"], [\"95\u00f78=11, 7\", \"58\u00f77=8, 2\"], [\"81\u00f73=27, 0\", \"52\u00f75=10, 2\"], [\"99\u00f74=24, 3\", \"70\u00f78=8, 6\"], [\"69\u00f73=23, 0\", \"69\u00f79=7, 6\"], [\"59\u00f76=9, 5\", \"17\u00f73=5, 2\"], [\"45\u00f72=22, 1\", \"30\u00f72=15, 0\"], [\"10\u00f75=2, 0\", \"34\u00f79=3, 7\"], [\"27\u00f76=4, 3\", \"98\u00f75=19, 3\"]];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet pairIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && pairIndex < pairs.length; i++) {\n  const para = paragraphs.items[i];\n  const current = para.text;\n  // Table rows added purely for spacing are empty paragraphs; skip them.\n  if (current === \"\") continue;\n\n  const [expectedOld, newText] = pairs[pairIndex];\n  if (current !== expectedOld) {\n    throw new Error(\n      `Unexpected paragraph text at index ${i} (pair ${pairIndex}): ` +\n      `expected ${JSON.stringify(expectedOld)}, found ${JSON.stringify(current)}`\n    );\n  }\n\n  para.getRange().insertText(newText, \"Replace\");\n  pairIndex++;\n}\n\nawait context.sync();\n\nif (pairIndex !== pairs.length) {\n  throw new Error(`Only replaced ${pairIndex} of ${pairs.length} expected paragraphs`);\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The document body is: one centered date paragraph, followed by a table of\n# two-digit division problems (5 columns x many rows, with several blank rows\n# used purely for spacing). Every non-blank paragraph changes text, in document\n# order, from an old value to a new value. $d.Paragraphs walks the whole story\n# (including paragraphs nested in table cells) in document order, so we can\n# replace them positionally while skipping the blank spacer paragraphs.\n\n$oldValues = @(\n    \"2023-08-02 Wednesday\",\n    \"89\u00f74=22, 1\",\n    \"89\u00f77=12, 5\",\n    \"36\u00f72=18, 0\",\n    \"28\u00f75=5, 3\",\n    \"86\u00f73=28, 2\",\n    \"99\u00f79=11, 0\",\n    \"38\u00f75=7, 3\",\n    \"98\u00f78=12, 2\",\n    \"82\u00f79=9, 1\",\n    \"80\u00f75=16, 0\",\n    \"83\u00f74=20, 3\",\n    \"79\u00f72=39, 1\",\n    \"18\u00f73=6, 0\",\n    \"89\u00f77=12, 5\",\n    \"93\u00f73=31, 0\",\n    \"64\u00f79=7, 1\",\n    \"83\u00f74=20, 3\",\n    \"95\u00f78=11, 7\",\n    \"81\u00f73=27, 0\",\n    \"99\u00f74=24, 3\",\n    \"69\u00f73=23, 0\",\n    \"59\u00f76=9, 5\",\n    \"45\u00f72=22, 1\",\n    \"10\u00f75=2, 0\",\n    \"27\u00f76=4, 3\",\n)\n$newValues = @(\n    \"2023-08-03 Thursday\",\n    \"87\u00f73=29, 0\",\n    \"12\u00f78=1, 4\",\n    \"97\u00f72=48, 1\",\n    \"72\u00f79=8, 0\",\n    \"38\u00f73=12, 2\",\n    \"16\u00f76=2, 4\",\n    \"62\u00f75=12, 2\",\n    \"60\u00f75=12, 0\",\n    \"71\u00f74=17, 3\",\n    \"12\u00f75=2, 2\",\n    \"27\u00f72=13, 1\",\n    \"89\u00f75=17, 4\",\n    \"51\u00f79=5, 6\",\n    \"75\u00f76=12, 3\",\n    \"24\u00f75=4, 4\",\n    \"93\u00f74=23, 1\",\n    \"32\u00f77=4, 4\",\n    \"58\u00f77=8, 2\",\n    \"52\u00f75=10, 2\",\n    \"70\u00f78=8, 6\",\n    \"69\u00f79=7, 6\",\n    \"17\u00f73=5, 2\",\n    \"30\u00f72=15, 0\",\n    \"34\u00f79=3, 7\",\n    \"98\u00f75=19, 3\",\n)\n\n$pairIndex = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count -and $pairIndex -lt $oldValues.Length; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $current = $r.Text\n    # Strip the trailing paragraph mark / cell mark before comparing.\n    $trimmed = $current.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"\") { continue }\n\n    $expectedOld = $oldValues[$pairIndex]\n    if ($trimmed -ne $expectedOld) {\n        throw \"Unexpected paragraph text at index $i (pair $pairIndex): expected [$expectedOld], found [$trimmed]\"\n    }\n\n    $r.Text = $newValues[$pairIndex]\n    $pairIndex++\n}\n\nif ($pairIndex -ne $oldValues.Length) {\n    throw \"Only replaced $pairIndex of $($oldValues.Length) expected paragraphs\"\n}\n\n"}
